# Reposition/resize the BEFORE/AFTER labels and the two comparison pictures
# on slide 1 (moving the picture split from vertical side-by-side to a
# stacked layout and nudging the labels to sit above the (now smaller)
# pictures).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Shape 1: "TextBox 3" (BEFORE label)
$shp = $s.Shapes.Item(1)
$shp.Left = 203.7276382446756
$shp.Top = 399.4636077881441

# Shape 2: "TextBox 4" (AFTER label)
$shp = $s.Shapes.Item(2)
$shp.Left = 698.8132083464567
$shp.Top = 399.4636077881441

# Shape 3: "Picture 1"
$shp = $s.Shapes.Item(3)
$shp.Left = 0.00007874015748031496
$shp.Top = 198.8571701050307
$shp.Width = 488.5713348388772
$shp.Height = 188.57141876223463

# Shape 4: "Picture 2"
$shp = $s.Shapes.Item(4)
$shp.Left = 488.57142639163465
$shp.Top = 198.8571701050307
$shp.Width = 471.42857360846534
$shp.Height = 188.57141876223463
